$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values look numeric need to be forced to Text format
# so Excel keeps them as strings (matching the original inline-string cells)
# instead of silently converting them to numbers.
$textForceCells = @('D4', 'D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D17', 'D18', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D41', 'D42', 'D43', 'D45', 'D47', 'D48', 'D50', 'D51')
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '29.035.00'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '1.828.70'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('D4').Value = '0.9985'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '240.71'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('D6').Value = '0.6201'
$ws.Range('E6').Value = '  -6.70%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = '0.07500'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '44.48'
$ws.Range('E9').Value = '  +6.14%  '
$ws.Range('D10').Value = '0.2911'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('D11').Value = '22.67'
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('D12').Value = '0.07618'
$ws.Range('E12').Value = '  -1.79%  '
$ws.Range('D13').Value = '1.826.65'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').Value = '4.948'
$ws.Range('E14').Value = '  -0.87%  '
$ws.Range('D15').Value = '0.6619'
$ws.Range('D16').Value = '81.96'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '0.000009138'
$ws.Range('E17').Value = '  +9.24%  '
$ws.Range('D18').Value = '5.969'
$ws.Range('E18').Value = '  -2.26%  '
$ws.Range('D19').Value = '29.030.30'
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').Value = '2.076.03'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').Value = '224.40'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').Value = '12.32'
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '7.166'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '159.42'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '8.401'
$ws.Range('E27').Value = '  -2.50%  '
$ws.Range('D28').Value = '0.1352'
$ws.Range('E28').Value = '  -4.38%  '
$ws.Range('D29').Value = '17.81'
$ws.Range('E29').Value = '  -0.98%  '
$ws.Range('D30').Value = '1.493'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('D31').Value = '4.029'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('D33').Value = '4.041'
$ws.Range('E33').Value = '  -1.72%  '
$ws.Range('D34').Value = '0.05208'
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').Value = '1.829'
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('D36').Value = '1.149'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('D37').Value = '0.7324'
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('D38').Value = '2.643'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').Value = '1.273.85'
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').Value = '0.01782'
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('D42').Value = '6.315'
$ws.Range('E42').Value = '  +7.37%  '
$ws.Range('D43').Value = '0.8937'
$ws.Range('E43').Value = '  -4.12%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = '101.75'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('D46').Value = '1.975.06'
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('D47').Value = '0.5121'
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('D48').Value = '63.29'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('D50').Value = '0.3956'
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('D51').Value = '1.676'
$ws.Range('E51').Value = '  -4.82%  '
